# Thesis binary-classification results workbook — add the "Top 5 models with
# lemmatization" block for the Lucene sheet (two new configurations, 5 models
# each), matching the commit "Top 5 models results with lemmatization for Lucene".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lucene")

# --- Widen the "Configuration" column (B) to fit the new, longer labels ----
$ws.Columns.Item(2).ColumnWidth = 38.0663265306122

# --- New data rows 57-66 ----------------------------------------------------
# Columns: A=Model, B=Configuration, C=F1, D=Accuracy, E=Precision, F=Recall

$config1 = "CV + tfidf + ngram(3) + stopwords + lem + RUS + SMOTE"
$config2 = "CV + tfidf + ngram(3) + stopwords + lem + SMOTE + RUS"

$rows = @(
    @("Logistic Regression",     $config1, 86.16, 79.96, 85.47, 87.1),
    @("Multinomial Naive Bayes", $config1, 84.34, 78.49, 88.01, 81.21),
    @("Support Vector Machines", $config1, 87.49, 79.88, 78.7,  98.86),
    @("Decision Tree",           $config1, 69.88, 63.89, 88.26, 58.07),
    @("Random Forest",           $config1, 70.94, 64.42, 86.66, 60.54),
    @("Logistic Regression",     $config2, 87.77, 82,    85.62, 90.21),
    @("Multinomial Naive Bayes", $config2, 87.85, 82.17, 85.64, 90.32),
    @("Support Vector Machines", $config2, 87.23, 79.31, 78.08, 99.31),
    @("Decision Tree",           $config2, 82,    74.69, 84.14, 80.22),
    @("Random Forest",           $config2, 79.48, 72.55, 86.3,  74.14)
)

$r = 57
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Move the selection/scroll position back to the top of the sheet and
# select the last filled cell, like the source file does after the edit. ----
$ws.Range("F66").Select()

# --- Cosmetic: widen the tab-bar a touch (tabRatio 986 -> 990) -------------
$excel.ActiveWindow.TabRatio = 990
